$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 68.443746
$ws.Range("H2").Value = 205.331238
$ws.Range("I2").Value = 0.1596169534001499
$ws.Range("J2").Value = 0.1596169534001499
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.604331
$ws.Range("N2").Value = 40.812993
$ws.Range("O2").Value = 0.8107276168878804
$ws.Range("P2").Value = 0.8107276168878805
$ws.Range("Q2").Value = 931.131375463926
$ws.Range("R2").Value = 8380.182379175334
$ws.Range("S2").Value = 0.1294058722450074
$ws.Range("T2").Value = 0.1294058722450074

# Row 3
$ws.Range("G3").Value = 68.443746
$ws.Range("H3").Value = 205.331238
$ws.Range("I3").Value = 0.1596169534001499
$ws.Range("J3").Value = 0.1596169534001499
$ws.Range("O3").Value = 0.06327311690486458
$ws.Range("P3").Value = 0.06327311690486459
$ws.Range("Q3").Value = 72.670010428008
$ws.Range("R3").Value = 654.030093852072
$ws.Range("S3").Value = 0.01009946215248601
$ws.Range("T3").Value = 0.01009946215248601

# Row 4
$ws.Range("G4").Value = 68.443746
$ws.Range("H4").Value = 205.331238
$ws.Range("I4").Value = 0.1596169534001499
$ws.Range("J4").Value = 0.1596169534001499
$ws.Range("M4").Value = 1.995771333333333
$ws.Range("N4").Value = 5.987314
$ws.Range("O4").Value = 0.1189346934389115
$ws.Range("P4").Value = 0.1189346934389116
$ws.Range("Q4").Value = 136.598066212748
$ws.Range("R4").Value = 1229.382595914732
$ws.Range("S4").Value = 0.01898399342029986
$ws.Range("T4").Value = 0.01898399342029986

# Row 5
$ws.Range("G5").Value = 68.443746
$ws.Range("H5").Value = 205.331238
$ws.Range("I5").Value = 0.1596169534001499
$ws.Range("J5").Value = 0.1596169534001499
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1185463333333334
$ws.Range("N5").Value = 0.355639
$ws.Range("O5").Value = 0.007064572768343379
$ws.Range("P5").Value = 0.007064572768343379
$ws.Range("Q5").Value = 8.113755127898001
$ws.Range("R5").Value = 73.02379615108201
$ws.Range("S5").Value = 0.001127625582356633
$ws.Range("T5").Value = 0.001127625582356633

# Row 6
$ws.Range("I6").Value = 0.4159650732941736
$ws.Range("J6").Value = 0.4159650732941736
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.604331
$ws.Range("N6").Value = 40.812993
$ws.Range("O6").Value = 0.8107276168878804
$ws.Range("P6").Value = 0.8107276168878805
$ws.Range("Q6").Value = 2426.547572740434
$ws.Range("R6").Value = 21838.92815466391
$ws.Range("S6").Value = 0.3372343725803779
$ws.Range("T6").Value = 0.337234372580378

# Row 7
$ws.Range("I7").Value = 0.4159650732941736
$ws.Range("J7").Value = 0.4159650732941736
$ws.Range("O7").Value = 0.06327311690486458
$ws.Range("P7").Value = 0.06327311690486459
$ws.Range("S7").Value = 0.02631940671088281
$ws.Range("T7").Value = 0.02631940671088282

# Row 8
$ws.Range("I8").Value = 0.4159650732941736
$ws.Range("J8").Value = 0.4159650732941736
$ws.Range("M8").Value = 1.995771333333333
$ws.Range("N8").Value = 5.987314
$ws.Range("O8").Value = 0.1189346934389115
$ws.Range("P8").Value = 0.1189346934389116
$ws.Range("Q8").Value = 355.9773784278653
$ws.Range("R8").Value = 3203.796405850788
$ws.Range("S8").Value = 0.04947267847353691
$ws.Range("T8").Value = 0.04947267847353692

# Row 9
$ws.Range("I9").Value = 0.4159650732941736
$ws.Range("J9").Value = 0.4159650732941736
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1185463333333334
$ws.Range("N9").Value = 0.355639
$ws.Range("O9").Value = 0.007064572768343379
$ws.Range("P9").Value = 0.007064572768343379
$ws.Range("Q9").Value = 21.14461324171534
$ws.Range("R9").Value = 190.301519175438
$ws.Range("S9").Value = 0.002938615529375977
$ws.Range("T9").Value = 0.002938615529375977

# Row 10
$ws.Range("G10").Value = 88.88346833333333
$ws.Range("H10").Value = 266.650405
$ws.Range("I10").Value = 0.2072842188241036
$ws.Range("J10").Value = 0.2072842188241036
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.604331
$ws.Range("N10").Value = 40.812993
$ws.Range("O10").Value = 0.8107276168878804
$ws.Range("P10").Value = 0.8107276168878805
$ws.Range("Q10").Value = 1209.200123634685
$ws.Range("R10").Value = 10882.80111271216
$ws.Range("S10").Value = 0.1680510407457314
$ws.Range("T10").Value = 0.1680510407457314

# Row 11
$ws.Range("G11").Value = 88.88346833333333
$ws.Range("H11").Value = 266.650405
$ws.Range("I11").Value = 0.2072842188241036
$ws.Range("J11").Value = 0.2072842188241036
$ws.Range("O11").Value = 0.06327311690486458
$ws.Range("P11").Value = 0.06327311690486459
$ws.Range("Q11").Value = 94.37184473597999
$ws.Range("R11").Value = 849.3466026238199
$ws.Range("S11").Value = 0.01311551861019104
$ws.Range("T11").Value = 0.01311551861019104

# Row 12
$ws.Range("G12").Value = 88.88346833333333
$ws.Range("H12").Value = 266.650405
$ws.Range("I12").Value = 0.2072842188241036
$ws.Range("J12").Value = 0.2072842188241036
$ws.Range("M12").Value = 1.995771333333333
$ws.Range("N12").Value = 5.987314
$ws.Range("O12").Value = 0.1189346934389115
$ws.Range("P12").Value = 0.1189346934389116
$ws.Range("Q12").Value = 177.3910781069077
$ws.Range("R12").Value = 1596.51970296217
$ws.Range("S12").Value = 0.02465328502056901
$ws.Range("T12").Value = 0.02465328502056902

# Row 13
$ws.Range("G13").Value = 88.88346833333333
$ws.Range("H13").Value = 266.650405
$ws.Range("I13").Value = 0.2072842188241036
$ws.Range("J13").Value = 0.2072842188241036
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1185463333333334
$ws.Range("N13").Value = 0.355639
$ws.Range("O13").Value = 0.007064572768343379
$ws.Range("P13").Value = 0.007064572768343379
$ws.Range("Q13").Value = 10.53680926486611
$ws.Range("R13").Value = 94.831283383795
$ws.Range("S13").Value = 0.001464374447612092
$ws.Range("T13").Value = 0.001464374447612092

# Row 14
$ws.Range("G14").Value = 93.106949
$ws.Range("H14").Value = 279.320847
$ws.Range("I14").Value = 0.2171337544815728
$ws.Range("J14").Value = 0.2171337544815728
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 13.604331
$ws.Range("N14").Value = 40.812993
$ws.Range("O14").Value = 0.8107276168878804
$ws.Range("P14").Value = 0.8107276168878805
$ws.Range("Q14").Value = 1266.657752596119
$ws.Range("R14").Value = 11399.91977336507
$ws.Range("S14").Value = 0.1760363313167637
$ws.Range("T14").Value = 0.1760363313167637

# Row 15
$ws.Range("G15").Value = 93.106949
$ws.Range("H15").Value = 279.320847
$ws.Range("I15").Value = 0.2171337544815728
$ws.Range("J15").Value = 0.2171337544815728
$ws.Range("O15").Value = 0.06327311690486458
$ws.Range("P15").Value = 0.06327311690486459
$ws.Range("Q15").Value = 98.85611688685199
$ws.Range("R15").Value = 889.705051981668
$ws.Range("S15").Value = 0.01373872943130472
$ws.Range("T15").Value = 0.01373872943130472

# Row 16
$ws.Range("G16").Value = 93.106949
$ws.Range("H16").Value = 279.320847
$ws.Range("I16").Value = 0.2171337544815728
$ws.Range("J16").Value = 0.2171337544815728
$ws.Range("M16").Value = 1.995771333333333
$ws.Range("N16").Value = 5.987314
$ws.Range("O16").Value = 0.1189346934389115
$ws.Range("P16").Value = 0.1189346934389116
$ws.Range("Q16").Value = 185.8201797483287
$ws.Range("R16").Value = 1672.381617734958
$ws.Range("S16").Value = 0.02582473652450575
$ws.Range("T16").Value = 0.02582473652450576

# Row 17
$ws.Range("G17").Value = 93.106949
$ws.Range("H17").Value = 279.320847
$ws.Range("I17").Value = 0.2171337544815728
$ws.Range("J17").Value = 0.2171337544815728
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1185463333333334
$ws.Range("N17").Value = 0.355639
$ws.Range("O17").Value = 0.007064572768343379
$ws.Range("P17").Value = 0.007064572768343379
$ws.Range("Q17").Value = 11.03748741180367
$ws.Range("R17").Value = 99.33738670623302
$ws.Range("S17").Value = 0.001533957208998677
$ws.Range("T17").Value = 0.001533957208998677
